# Fix error of turm classes being only one
# Clear out the duplicated/erroneous class entries in the MCT-3A schedule,
# replacing them with "-" to indicate no class in that slot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "-"
$ws.Range("C3").Value = "-"
$ws.Range("C7").Value = "-"
$ws.Range("D10").Value = "-"
